$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the changed "Price" (D) cells keep their original literal-text
# representation (e.g. trailing zeros, grouped-looking numbers such as
# "30.195.13") instead of being auto-converted to a number by Excel's normal
# value-parsing heuristics. Column E values are already non-numeric text
# (they carry padding spaces + a "%" sign) so a plain .Value assignment is safe.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.195.13'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.911.65'
$ws.Range("E3").Value = '  -0.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.13'
$ws.Range("E5").Value = '  -2.91%  '

$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5064'
$ws.Range("E7").Value = '  -2.93%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4040'
$ws.Range("E8").Value = '  -1.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08285'
$ws.Range("E9").Value = '  -2.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.111'
$ws.Range("E10").Value = '  -1.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.28'
$ws.Range("E11").Value = '  -1.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.99'
$ws.Range("E12").Value = '  +3.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.417'
$ws.Range("E13").Value = '  -0.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.894.66'
$ws.Range("E14").Value = '  -0.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.327'
$ws.Range("E15").Value = '  -1.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.71'
$ws.Range("E17").Value = '  -2.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001099'
$ws.Range("E18").Value = '  -1.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06471'
$ws.Range("E19").Value = '  -3.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.55'
$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.953'
$ws.Range("E22").Value = '  -1.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.213.70'
$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.31'
$ws.Range("E24").Value = '  -0.63%  '

$ws.Range("E25").Value = '  -0.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '22.31'
$ws.Range("E26").Value = '  +5.33%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.119.75'
$ws.Range("E27").Value = '  -0.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.86'
$ws.Range("E28").Value = '  -0.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.368'
$ws.Range("E29").Value = '  -2.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.91'
$ws.Range("E30").Value = '  +0.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.124'
$ws.Range("E31").Value = '  +2.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1044'
$ws.Range("E32").Value = '  -2.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.000'
$ws.Range("E33").Value = '  -0.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.813'
$ws.Range("E34").Value = '  +5.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02462'
$ws.Range("E35").Value = '  -1.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.404'
$ws.Range("E36").Value = '  +4.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06465'
$ws.Range("E37").Value = '  -1.96%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2161'
$ws.Range("E38").Value = '  -2.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.741'
$ws.Range("E39").Value = '  -1.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.191'
$ws.Range("E40").Value = '  -3.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6407'
$ws.Range("E41").Value = '  -2.10%  '

$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.38'
$ws.Range("E42").Value = '  -4.25%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.217'
$ws.Range("E43").Value = '  -2.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9996'
$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.198'
$ws.Range("E45").Value = '  +4.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.23'
$ws.Range("E46").Value = '  -0.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6001'
$ws.Range("E47").Value = '  -2.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.649'
$ws.Range("E48").Value = '  -2.69%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.79'
$ws.Range("E49").Value = '  -0.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.218'
$ws.Range("E50").Value = '  -2.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.03'
$ws.Range("E51").Value = '  -0.93%  '
